$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text (matches original inlineStr / shared-string text cells)
# by forcing a text number format before assigning values that look numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.779.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.86%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.60%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.75%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.681'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.88%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.88%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.67'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.36%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '51.05'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0734'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.87%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '12.86'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.145.51'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.713'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.78%  '

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.86'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.871.01'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '34.789.39'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.54'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0816'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '242.60'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.64'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.72%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.81%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.83'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.32'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.03%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.99%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.128.43'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.20'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0575'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.14'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.45%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.824'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -9.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.97'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.51'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -23.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '97.08'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.90'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.98%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0665'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.06%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.75%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.281.64'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.31'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0800'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +9.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.72'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.94'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.87%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.37%  '
